# "Actualizar" automation run: refreshes the availability-check timestamp
# for the most-recent 14-row batch (rows 212-225) and appends a brand new
# batch of 14 rows (226-239) with the same Nombre/URL/Disponibilidad cycle,
# each carrying its own hyperlink in column B.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- 1. Refresh the timestamp on the previous batch (rows 212-225) ----
$oldTimestamp = 44231.98718089121
for ($r = 212; $r -le 225; $r++) {
    $ws.Range("D$r").Value = $oldTimestamp
}

# ---- 2. Append the new batch (rows 226-239) ----
$names = @("Odoo","Blackbox","PowerBI","Dropbox","Odoo","GEE","UtilidadesOdoo","Filtros Dashboard","MapStore","GeoServer","Tomcat","Shiny","Github","EZ Exporter")
$urls  = @(
    "https://www.dataintelligence-group.com/",
    "https://serviciodashboard.azurewebsites.net/",
    "https://powerbi.microsoft.com/es-es/",
    "https://www.dropbox.com/",
    "https://dataintelligence.store/",
    "https://app-data-i.users.earthengine.app/",
    "https://odooutil.azurewebsites.net/",
    "https://filtradordashboard.azurewebsites.net/",
    "https://ide.dataintelligence-group.com/mapstore/",
    "https://ide.dataintelligence-group.com/geoserver/web/?0",
    "https://ide.dataintelligence-group.com/",
    "https://rpubs.com/dataintelligence/",
    "https://github.com/Sud-Austral/",
    "https://ezexporter.highviewapps.com/exports/export-profile/"
)
# Only the MapStore entry (index 8) carries an in-page fragment ("#/"),
# stored as SubAddress/location rather than baked into the Target URL.
$subAddresses = @("","","","","","","","","/","","","","","")

$newTimestamp = 44232.00822890468
$startRow = 226

for ($i = 0; $i -lt $names.Length; $i++) {
    $r = $startRow + $i

    $ws.Range("A$r").Value = $names[$i]
    $ws.Range("C$r").Value = "Disponible"

    $ws.Range("D$r").Value = $newTimestamp
    $ws.Range("D$r").NumberFormat = "YYYY-MM-DD HH:MM:SS"

    if ($subAddresses[$i] -ne "") {
        $ws.Hyperlinks.Add($ws.Range("B$r"), $urls[$i], $subAddresses[$i])
        $ws.Range("B$r").Value = $urls[$i] + "#" + $subAddresses[$i]
    } else {
        $ws.Hyperlinks.Add($ws.Range("B$r"), $urls[$i])
        $ws.Range("B$r").Value = $urls[$i]
    }
    $ws.Range("B$r").Style = "Hyperlink"
}
